$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-23 00:48:10"
$ws.Rows.Item(3).Delete()

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-23 00:48:06"
$ws.Rows.Item(3).Delete()

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-23 00:48:10"
$ws.Rows.Item(3).Delete()
